# Weekly update: add the latest "Fruta, Vega Modelo de Temuco - Membrillo"
# price record as a new first data row, pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data block (row 78) - this shifts all
# existing data rows (78-143) down by one (to 79-144) and keeps their
# values/styles intact.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44658
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100104
$ws.Range("H78").Value = "Frutos de pepita"
$ws.Range("I78").Value = 100104003
$ws.Range("J78").Value = "Membrillo"
$ws.Range("K78").Value = "Champion"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 200
$ws.Range("N78").Value = 12000
$ws.Range("O78").Value = 12000
$ws.Range("P78").Value = 12000
$ws.Range("Q78").Value = "$/bandeja 18 kilos granel"
$ws.Range("R78").Value = "Región de O'Higgins"
$ws.Range("S78").Value = 667
$ws.Range("T78").Value = 18
